$wb = $excel.ActiveWorkbook

# --- Rename worksheets (Weekly -> Monthly) ---
$wsFirst = $wb.Worksheets.Item("first")
$wsGeneral = $wb.Worksheets.Item("GeneralTaxRateWeekly")
$wsProcess = $wb.Worksheets.Item("ProcessPayrollForWeeklyTax")
$wsReports = $wb.Worksheets.Item("TestReports")

$wsGeneral.Name = "GeneralTaxRateMonthly"
$wsProcess.Name = "ProcessPayrollForMonthlyTax"

# --- Update cell contents referencing the old sheet names / labels ---
$wsFirst.Range("A3").Value = "GeneralTaxRateMonthly"
$wsFirst.Range("A4").Value = "ProcessPayrollForMonthlyTax"

$wsGeneral.Range("A2").Value = "DO NOT TOUCH AUTOMATION EMP 105"
$wsProcess.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"
$wsReports.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"

# --- Update view state: active tab moves to GeneralTaxRateMonthly (index 2) ---
$wsGeneral.Activate()
$wsGeneral.Range("E11").Select()

$wsFirst.Activate()
$wsFirst.Range("F5").Select()

$wsProcess.Activate()
$wsProcess.Range("D5").Select()

$wsReports.Activate()
$wsReports.Range("F5").Select()

# Final active sheet / tab should be GeneralTaxRateMonthly
$wsGeneral.Activate()
